$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the duplicate "Contact" / "No display for ContactDetail" row (row 11);
# this shifts Description/Purpose/Copyright/Immutable rows up by one.
$ws.Rows.Item(11).Delete()

# Update Version
$ws.Range("B3").Value = "6.0.0"

# Update Date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Set Publisher value
$ws.Range("B9").Value = "Alvearie Team"

# Replace the remaining Contact row with Jurisdiction info
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
